$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Insert a new paragraph right after "Das Layout/Design zur Webanwendung
# liegt im Anhang bei." containing the new text about the login-protected
# web interface, including the redirect sentence. Word's "_GoBack" bookmark
# (which marks the last edit position) ends up splitting the final run into
# "...Login-Seite um" | "geleitet.", matching the target XML exactly.
# ---------------------------------------------------------------------------

$anchor = $d.Content
$found = $anchor.Find.Execute("Das Layout/Design zur Webanwendung liegt im Anhang bei.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

$newPara = $d.Paragraphs(4)
$insertionPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertionPoint.InsertAfter("Die Weboberfläche steht ausschließlich eingeloggten Nutzern vor.")

$newPara = $d.Paragraphs(4)
$tail = [char]11 + "Nicht angemeldete Anwender werden grundsätzlich bei Zugriff auf die Seite auf die Login-Seite umgeleitet."
$tailPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$tailPoint.InsertAfter($tail)

# Re-place the singleton "_GoBack" bookmark exactly between "...Login-Seite um"
# and "geleitet." -- this both moves it out of its old position (elsewhere in
# the document) and splits the run the way a live edit naturally would.
$newPara = $d.Paragraphs(4)
$lastWord = "geleitet."
$bookmarkPos = $newPara.Range.End - 1 - $lastWord.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

Write-Output "Inserted paragraph text: [$($d.Paragraphs(4).Range.Text)]"
